# Updating filtered feeds from workflow
# Appends a new feed entry (link / keywords / title) as row 3 of the
# "Filtered Feeds" sheet, mirroring the format of the existing row 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newLink    = "https://www.360dx.com/cancer/agilent-gains-ivdr-certification-expanded-use-cdx-assay-keytruda"
$newKeyword = "CDx"
$newTitle   = "Agilent Gains IVDR Certification for Expanded Use of CDx Assay for Keytruda"

# Write the link text first, then turn it into a real hyperlink (same
# pattern Excel uses for A2), and finally reapply the built-in "Hyperlink"
# cell style so A3 matches A2's formatting exactly instead of picking up a
# freshly-minted (but equivalent) style record.
$ws.Range("A3").Value = $newLink
$ws.Hyperlinks.Add($ws.Range("A3"), $newLink)
$ws.Range("A3").Style = "Hyperlink"

$ws.Range("B3").Value = $newKeyword
$ws.Range("C3").Value = $newTitle
